# updated INTC, PLTR, created APP
# Bump the "last updated" dates on the Hardware & Semis (D11), Software (D12),
# Services (D8), and Gaming (D3) sector screens, plus the Hedge Funds date (C19),
# to 12/9/2024 (Excel serial 45635).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Range("D3").Value = 45635
$ws.Range("D8").Value = 45635
$ws.Range("D11").Value = 45635
$ws.Range("D12").Value = 45635
$ws.Range("C19").Value = 45635

$ws.Range("D19").Select()
